$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: update existing cell values (recalculated statistics) ---
$updates = @(
    @{row=3; col=7; val=0.0092765977089689},
    @{row=4; col=7; val=0.0092765977089689},
    @{row=9; col=7; val=0.0144894851136709},
    @{row=10; col=7; val=0.0144894851136709},
    @{row=13; col=6; val=0.6648500000000001},
    @{row=13; col=7; val=0.664203571428571},
    @{row=13; col=9; val=1.03195},
    @{row=13; col=14; val=1.00234},
    @{row=14; col=6; val=0.6648500000000001},
    @{row=14; col=7; val=0.664203571428571},
    @{row=14; col=9; val=1.03195},
    @{row=14; col=14; val=1.00234},
    @{row=20; col=7; val=0.0085936406195116},
    @{row=21; col=7; val=0.0085936406195116},
    @{row=26; col=7; val=0.0129109175921514},
    @{row=26; col=12; val=0.00234},
    @{row=27; col=7; val=0.0129109175921514},
    @{row=27; col=12; val=0.00234},
    @{row=30; col=6; val=0.6756},
    @{row=30; col=7; val=0.689842372881356},
    @{row=30; col=9; val=1.02812},
    @{row=30; col=13; val=0.97827},
    @{row=30; col=14; val=1.00608},
    @{row=31; col=6; val=0.6756},
    @{row=31; col=7; val=0.689842372881356},
    @{row=31; col=9; val=1.02812},
    @{row=31; col=13; val=0.97827},
    @{row=31; col=14; val=1.00608},
    @{row=37; col=7; val=0.008809839482932199},
    @{row=38; col=7; val=0.008809839482932199},
    @{row=43; col=7; val=0.009915228167610301},
    @{row=43; col=12; val=0.00172},
    @{row=44; col=7; val=0.009915228167610301},
    @{row=44; col=12; val=0.00172},
    @{row=47; col=7; val=0.703161666666667},
    @{row=47; col=12; val=0.45045},
    @{row=47; col=14; val=1.04338},
    @{row=48; col=7; val=0.703161666666667},
    @{row=48; col=12; val=0.45045},
    @{row=48; col=14; val=1.04338},
    @{row=54; col=7; val=0.008894456318692699},
    @{row=55; col=7; val=0.008894456318692699},
    @{row=60; col=6; val=0.00236},
    @{row=60; col=7; val=0.0068856929712089},
    @{row=60; col=12; val=0.00172},
    @{row=61; col=6; val=0.00236},
    @{row=61; col=7; val=0.0068856929712089},
    @{row=61; col=12; val=0.00172},
    @{row=64; col=7; val=0.7037949999999999},
    @{row=64; col=12; val=0.425},
    @{row=64; col=13; val=0.99013},
    @{row=65; col=7; val=0.7037949999999999},
    @{row=65; col=12; val=0.425},
    @{row=65; col=13; val=0.99013},
    @{row=71; col=7; val=0.009457356037389999},
    @{row=72; col=7; val=0.009457356037389999},
    @{row=77; col=6; val=0.00301},
    @{row=77; col=7; val=0.005855034920745},
    @{row=77; col=12; val=0.00236},
    @{row=78; col=6; val=0.00301},
    @{row=78; col=7; val=0.005855034920745},
    @{row=78; col=12; val=0.00236},
    @{row=81; col=7; val=0.706342372881356},
    @{row=81; col=12; val=0.42925},
    @{row=81; col=13; val=0.99304},
    @{row=82; col=7; val=0.706342372881356},
    @{row=82; col=12; val=0.42925},
    @{row=82; col=13; val=0.99304},
    @{row=87; col=7; val=1.40880238758664},
    @{row=87; col=13; val=2.247},
    @{row=88; col=7; val=0.0093050717523396},
    @{row=89; col=7; val=0.0093050717523396},
    @{row=94; col=7; val=0.0058250395613937},
    @{row=94; col=14; val=0.01337},
    @{row=95; col=7; val=0.0058250395613937},
    @{row=95; col=14; val=0.01337},
    @{row=98; col=7; val=0.710220689655172},
    @{row=99; col=7; val=0.710220689655172},
    @{row=104; col=7; val=1.29237809269575},
    @{row=104; col=13; val=1.9534},
    @{row=105; col=7; val=0.009054743185736},
    @{row=106; col=7; val=0.009054743185736},
    @{row=111; col=7; val=0.0071678054184639},
    @{row=111; col=12; val=0.007889999999999999},
    @{row=112; col=7; val=0.0071678054184639},
    @{row=112; col=12; val=0.007889999999999999},
    @{row=115; col=7; val=0.6824375},
    @{row=116; col=7; val=0.6824375},
    @{row=121; col=7; val=0.917859654925435},
    @{row=121; col=9; val=2.29697},
    @{row=121; col=14; val=1.908},
    @{row=122; col=7; val=0.008268335402038101},
    @{row=123; col=7; val=0.008268335402038101},
    @{row=128; col=6; val=0.00645},
    @{row=128; col=7; val=0.008262725256161301},
    @{row=128; col=12; val=0.007889999999999999},
    @{row=129; col=6; val=0.00645},
    @{row=129; col=7; val=0.008262725256161301},
    @{row=129; col=12; val=0.007889999999999999},
    @{row=138; col=7; val=0.008176995735571201},
    @{row=139; col=7; val=0.008176995735571201},
    @{row=140; col=7; val=586.463244705456},
    @{row=141; col=7; val=586.463244705456},
    @{row=142; col=7; val=586.463244705456},
    @{row=143; col=7; val=586.463244705456},
    @{row=144; col=6; val=0.00764},
    @{row=144; col=7; val=0.0089186653357449},
    @{row=144; col=13; val=0.01309},
    @{row=145; col=6; val=0.00764},
    @{row=145; col=7; val=0.0089186653357449},
    @{row=145; col=13; val=0.01309},
    @{row=148; col=12; val=0.51025},
    @{row=149; col=12; val=0.51025}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.row, $u.col).Value = $u.val
}

# --- Part 2: append new rows 154-169 (2019-2023 period results) ---
$newRows = @(
    @{r=154; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='DRP (95th Percentile)'}, @{c=3; v='A'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.0065}, @{c=7; v=0.0074106282381184}, @{c=8; v=0.019}, @{c=9; v=0.0141}, @{c=12; v=0.005}, @{c=13; v=0.011}, @{c=14; v=0.013}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='mg/L'})},
    @{r=155; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='DRP (Median)'}, @{c=3; v='B'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.0065}, @{c=7; v=0.0074106282381184}, @{c=8; v=0.019}, @{c=9; v=0.0141}, @{c=12; v=0.005}, @{c=13; v=0.011}, @{c=14; v=0.013}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='mg/L'})},
    @{r=156; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='E coli (>260)'}, @{c=3; v='C'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=135}, @{c=7; v=338.51303083288}, @{c=8; v=2200}, @{c=9; v=1500}, @{c=10; v=22.9166666666667}, @{c=11; v=31.25}, @{c=12; v=52}, @{c=13; v=708.36}, @{c=14; v=1198}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='% exceedances over 260/100 mL'})},
    @{r=157; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='E coli (>540)'}, @{c=3; v='D'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=135}, @{c=7; v=338.51303083288}, @{c=8; v=2200}, @{c=9; v=1500}, @{c=10; v=22.9166666666667}, @{c=11; v=31.25}, @{c=12; v=52}, @{c=13; v=708.36}, @{c=14; v=1198}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='% exceedances over 540/100 mL'})},
    @{r=158; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='E coli (Median)'}, @{c=3; v='D'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=135}, @{c=7; v=338.51303083288}, @{c=8; v=2200}, @{c=9; v=1500}, @{c=10; v=22.9166666666667}, @{c=11; v=31.25}, @{c=12; v=52}, @{c=13; v=708.36}, @{c=14; v=1198}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='E. coli/100 mL'})},
    @{r=159; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='E coli (95th Percentile)'}, @{c=3; v='E'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=135}, @{c=7; v=338.51303083288}, @{c=8; v=2200}, @{c=9; v=1500}, @{c=10; v=22.9166666666667}, @{c=11; v=31.25}, @{c=12; v=52}, @{c=13; v=708.36}, @{c=14; v=1198}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='E. coli/100 mL'})},
    @{r=160; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Ammoniacal-N (95th Percentile)'}, @{c=3; v='A'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.007889999999999999}, @{c=7; v=0.009438736902649099}, @{c=8; v=0.0329027150070411}, @{c=9; v=0.02411}, @{c=12; v=0.009039999999999999}, @{c=13; v=0.01463}, @{c=14; v=0.02231}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='mg NH4-N/L'})},
    @{r=161; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Ammoniacal-N (Median)'}, @{c=3; v='A'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.007889999999999999}, @{c=7; v=0.009438736902649099}, @{c=8; v=0.0329027150070411}, @{c=9; v=0.02411}, @{c=12; v=0.009039999999999999}, @{c=13; v=0.01463}, @{c=14; v=0.02231}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='mg NH4-N/L'})},
    @{r=162; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Nitrate-N (95th Percentile)'}, @{c=3; v='A'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.631}, @{c=7; v=0.653}, @{c=8; v=1.13}, @{c=9; v=1.015}, @{c=12; v=0.504}, @{c=13; v=0.8801}, @{c=14; v=0.99116}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='mg NO3-N/L'})},
    @{r=163; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Nitrate-N (Median)'}, @{c=3; v='A'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.631}, @{c=7; v=0.653}, @{c=8; v=1.13}, @{c=9; v=1.015}, @{c=12; v=0.504}, @{c=13; v=0.8801}, @{c=14; v=0.99116}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='mg NO3-N/L'})},
    @{r=164; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Soluble Inorganic Nitrogen (95th Percentile)'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.647}, @{c=7; v=0.666958333333333}, @{c=8; v=1.159}, @{c=9; v=1.0459}, @{c=12; v=0.515}, @{c=13; v=0.88846}, @{c=14; v=1.0048}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='g/m3'})},
    @{r=165; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Soluble Inorganic Nitrogen (Median)'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.647}, @{c=7; v=0.666958333333333}, @{c=8; v=1.159}, @{c=9; v=1.0459}, @{c=12; v=0.515}, @{c=13; v=0.88846}, @{c=14; v=1.0048}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='g/m3'})},
    @{r=166; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Total Nitrogen (95th Percentile)'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.8}, @{c=7; v=0.837916666666667}, @{c=8; v=1.35}, @{c=9; v=1.269}, @{c=12; v=0.6899999999999999}, @{c=13; v=1.027}, @{c=14; v=1.1464}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='g/m3'})},
    @{r=167; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Total Nitrogen (Median)'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.8}, @{c=7; v=0.837916666666667}, @{c=8; v=1.35}, @{c=9; v=1.269}, @{c=12; v=0.6899999999999999}, @{c=13; v=1.027}, @{c=14; v=1.1464}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='g/m3'})},
    @{r=168; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Total Phosphorus (95th Percentile)'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.021}, @{c=7; v=0.0261875}, @{c=8; v=0.098}, @{c=9; v=0.065}, @{c=12; v=0.015}, @{c=13; v=0.04034}, @{c=14; v=0.0556}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='g/m3'})},
    @{r=169; cells=@(@{c=1; v='Tiraumea u/s Manawatu Confluence'}, @{c=2; v='Total Phosphorus (Median)'}, @{c=4; v='2019 - 2023'}, @{c=5; v='RepSite'}, @{c=6; v=0.021}, @{c=7; v=0.0261875}, @{c=8; v=0.098}, @{c=9; v=0.065}, @{c=12; v=0.015}, @{c=13; v=0.04034}, @{c=14; v=0.0556}, @{c=15; v=1845196}, @{c=16; v=5525095}, @{c=17; v='Tararua District'}, @{c=18; v='Manawatū'}, @{c=19; v='Tiraumea'}, @{c=20; v='Mana_7b'}, @{c=21; v='g/m3'})}
)

foreach ($rowDef in $newRows) {
    foreach ($cell in $rowDef.cells) {
        $ws.Cells.Item($rowDef.r, $cell.c).Value = $cell.v
    }
}
